# Update rows 121-126: mark as "VOLADA" with a last-modified date instead of "PENDIENTE"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 121; $r -le 126; $r++) {
    $ws.Cells.Item($r, 2).Value = "VOLADA"
    $ws.Cells.Item($r, 3).Value = 45751
}
